$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 'Alejandro E. Ulvert (3058421968) (D), Adi (3103038889) (D)'
$ws.Cells.Item(3, 3).Value = 'Alexander (9175433831) (D), Alejandro L (2676633214), Thor Waguespack (2404299065) (D)'
$ws.Cells.Item(6, 3).Value = 'Kamsi (6465933923), Alejandro Espinosa (7866702380) (D), Harry Corbin (5182223247) (D)'
$ws.Cells.Item(7, 3).Value = 'Ezana (2405594003), Paul (3106583636) (D), George Ryckman (9176015863) (D)'
$ws.Cells.Item(8, 3).Value = 'Jamari Pitchford (7738297627) (D), Edu (7862012521) (D), Noah Yaffe (6109995081) (D)'
$ws.Cells.Item(9, 3).Value = 'Adi (3103038889) (D), Harry Corbin (5182223247) (D)'
$ws.Cells.Item(10, 3).Value = 'Paul (3106583636) (D), Josh Greene (9176368454) (D), Alejandro E. Ulvert (3058421968) (D)'
$ws.Cells.Item(12, 3).Value = 'Ali Awada (6469771844), George Ryckman (9176015863) (D)'
$ws.Cells.Item(13, 3).Value = 'Alejandro Espinosa (7866702380) (D), Edu (7862012521) (D), Blake Steel (6109995081) (D)'
$ws.Cells.Item(14, 3).Value = 'Matheo (6462079196) (D), Jake Dieterich (9177968225), Noah Yaffe (6109995081) (D)'
$ws.Cells.Item(15, 3).Value = 'George Ryckman (9176015863) (D), Jake Dieterich (9177968225)'
$ws.Cells.Item(16, 3).Value = 'Alejandro L (2676633214), Alejandro Espinosa (7866702380) (D), Thor Waguespack (2404299065) (D)'
$ws.Cells.Item(17, 3).Value = 'Jaxon (4159881691) (D), Matheo (6462079196) (D)'
$ws.Cells.Item(18, 3).Value = 'Henry (6073399363), Edu (7862012521) (D), Blake Steel (6109995081) (D)'
$ws.Cells.Item(19, 3).Value = 'Gabe Heller (7187220895) (D), Kamsi (6465933923), Josh Greene (9176368454) (D)'
$ws.Cells.Item(20, 3).Value = 'Ben Kairouz (9179401653) (D), Alexander (9175433831) (D), Alejandro E. Ulvert (3058421968) (D)'
$ws.Cells.Item(21, 3).Value = 'Ezana (2405594003), Jamari Pitchford (7738297627) (D), Noah Yaffe (6109995081) (D)'
$ws.Cells.Item(22, 3).Value = 'Kamsi (6465933923), Blake Steel (6109995081) (D)'
$ws.Cells.Item(23, 3).Value = 'Gabe Heller (7187220895) (D), Josh Greene (9176368454) (D), Thor Waguespack (2404299065) (D)'
$ws.Cells.Item(25, 3).Value = 'Ben Kairouz (9179401653) (D), Ali Awada (6469771844)'
$ws.Cells.Item(26, 3).Value = 'Henry (6073399363), Thor Waguespack (2404299065) (D)'
$ws.Cells.Item(27, 3).Value = 'Jaxon (4159881691) (D), Gabe Heller (7187220895) (D), Alejandro Espinosa (7866702380) (D)'
$ws.Cells.Item(29, 3).Value = 'Ben Kairouz (9179401653) (D), Alejandro E. Ulvert (3058421968) (D), Harry Corbin (5182223247) (D)'
$ws.Cells.Item(30, 3).Value = 'Paul (3106583636) (D), Alexander (9175433831) (D), Alejandro L (2676633214)'
$ws.Cells.Item(31, 3).Value = 'Jamari Pitchford (7738297627) (D), Ali Awada (6469771844)'
$ws.Cells.Item(32, 3).Value = 'Jack Mogelof (9172163044) (D), Alejandro L (2676633214), George Ryckman (9176015863) (D)'
$ws.Cells.Item(33, 3).Value = 'Gabe Heller (7187220895) (D), Alexander (9175433831) (D), Edu (7862012521) (D)'
$ws.Cells.Item(34, 3).Value = 'Jaxon (4159881691) (D), Ezana (2405594003), Josh Greene (9176368454) (D)'
$ws.Cells.Item(35, 3).Value = 'Paul (3106583636) (D), Henry (6073399363), Adi (3103038889) (D)'
$ws.Cells.Item(36, 3).Value = 'Ezana (2405594003), Matheo (6462079196) (D)'
$ws.Cells.Item(37, 3).Value = 'Jaxon (4159881691) (D), Ben Kairouz (9179401653) (D), Kamsi (6465933923)'
$ws.Cells.Item(38, 3).Value = 'Paul (3106583636) (D), Ali Awada (6469771844), Harry Corbin (5182223247) (D)'
$ws.Cells.Item(39, 3).Value = 'Jack Mogelof (9172163044) (D), Gabe Heller (7187220895) (D), Jake Dieterich (9177968225)'
$ws.Cells.Item(40, 3).Value = 'Alejandro L (2676633214), Alejandro Espinosa (7866702380) (D), Noah Yaffe (6109995081) (D)'
$ws.Cells.Item(41, 3).Value = 'Alexander (9175433831) (D), Josh Greene (9176368454) (D), Henry (6073399363)'
